# Update the lattice-multiplication exercise table: every cell's problem
# (top line), the multiplicand split row, and the two partial-product rows
# get new values per the target revision. The "----" separator line is
# unchanged in every cell, so it is re-emitted as a literal.
#
# Each cell is addressed by (row, col) and the four lines that differ are
# supplied; InsertXML is used (rather than plain .Text assignment) so the
# xml:space="preserve" attribute on the space-padded lines round-trips
# exactly like the original markup.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cellData = @(
  @(1,1,"27 x 10","  1    0","2|    |","7|    |"),
  @(1,2,"86 x 47","  4    7","8|    |","6|    |"),
  @(1,3,"29 x 21","  2    1","2|    |","9|    |"),
  @(2,1,"93 x 54","  5    4","9|    |","3|    |"),
  @(2,2,"19 x 48","  4    8","1|    |","9|    |"),
  @(2,3,"56 x 46","  4    6","5|    |","6|    |"),
  @(3,1,"21 x 86","  8    6","2|    |","1|    |"),
  @(3,2,"52 x 89","  8    9","5|    |","2|    |"),
  @(3,3,"63 x 57","  5    7","6|    |","3|    |"),
  @(4,1,"11 x 48","  4    8","1|    |","1|    |"),
  @(4,2,"93 x 31","  3    1","9|    |","3|    |"),
  @(4,3,"95 x 35","  3    5","9|    |","5|    |"),
  @(5,1,"64 x 50","  5    0","6|    |","4|    |"),
  @(5,2,"86 x 53","  5    3","8|    |","6|    |"),
  @(5,3,"70 x 50","  5    0","7|    |","0|    |")
)

foreach ($row in $cellData) {
  $r = $row[0]
  $c = $row[1]
  $line1 = $row[2]
  $line2 = $row[3]
  $line4 = $row[4]
  $line5 = $row[5]

  $inner = '<w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr>' +
           '<w:t>{0}</w:t><w:br/>' +
           '<w:t xml:space="preserve">{1}</w:t><w:br/>' +
           '<w:t xml:space="preserve">  ----</w:t><w:br/>' +
           '<w:t>{2}</w:t><w:br/>' +
           '<w:t>{3}</w:t></w:r></w:p>'
  $inner = $inner -f $line1, $line2, $line4, $line5

  $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
         '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
         '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
         $inner +
         '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

  $cell = $t.Cell($r, $c)
  $cell.Range.InsertXML($xml)
}
